$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking price/volume cells are plain text in this sheet (e.g. "67.343.54"
# using "." as a thousands separator, or "  -1.49%  "). A handful of the new
# price values (like "1.00") would otherwise be auto-detected by Excel as
# numbers, so those are entered with a leading apostrophe to keep them text,
# exactly like a user typing '1.00 into the cell.
$ws.Range("D2").Value = '67.383.03'
$ws.Range("E2").Value = '  -1.40%  '
$ws.Range("D3").Value = '3.754.28'
$ws.Range("E3").Value = '  -2.14%  '
$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = "'595.09"
$ws.Range("E5").Value = '  -1.09%  '
$ws.Range("D6").Value = "'169.23"
$ws.Range("E6").Value = '  -0.33%  '
$ws.Range("D7").Value = '3.751.16'
$ws.Range("E7").Value = '  -2.27%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("E9").Value = '  -0.89%  '
$ws.Range("E10").Value = '  +0.12%  '
$ws.Range("D11").Value = "'6.49"
$ws.Range("E12").Value = '  -1.34%  '
$ws.Range("E13").Value = '  +3.77%  '
$ws.Range("D14").Value = "'36.43"
$ws.Range("E14").Value = '  -2.17%  '
$ws.Range("D15").Value = '4.385.42'
$ws.Range("E15").Value = '  -2.10%  '
$ws.Range("D16").Value = '3.753.94'
$ws.Range("E16").Value = '  -2.14%  '
$ws.Range("D17").Value = "'18.57"
$ws.Range("E17").Value = '  +0.56%  '
$ws.Range("D18").Value = '67.324.78'
$ws.Range("E18").Value = '  -1.63%  '
$ws.Range("E19").Value = '  -3.36%  '
$ws.Range("E20").Value = '  +0.74%  '
$ws.Range("D21").Value = "'10.53"
$ws.Range("E21").Value = '  -5.10%  '
$ws.Range("D22").Value = "'466.47"
$ws.Range("E22").Value = '  -0.65%  '
$ws.Range("D23").Value = "'0.719"
$ws.Range("E23").Value = '  -2.65%  '
$ws.Range("D25").Value = "'83.64"
$ws.Range("E25").Value = '  +0.45%  '
$ws.Range("E26").Value = '  -1.71%  '
$ws.Range("D27").Value = "'12.16"
$ws.Range("E27").Value = '  -0.74%  '
$ws.Range("E28").Value = '  +2.46%  '
$ws.Range("E29").Value = '  +0.02%  '
$ws.Range("E30").Value = '  -2.52%  '
$ws.Range("D31").Value = '3.902.06'
$ws.Range("E31").Value = '  -2.11%  '
$ws.Range("D32").Value = "'7.64"
$ws.Range("E32").Value = '  -0.91%  '
$ws.Range("D33").Value = "'30.43"
$ws.Range("E33").Value = '  -3.94%  '
$ws.Range("E34").Value = '  -4.10%  '
$ws.Range("D35").Value = "'9.11"
$ws.Range("E35").Value = '  -3.80%  '
$ws.Range("D36").Value = '3.717.93'
$ws.Range("E36").Value = '  -2.23%  '
$ws.Range("D37").Value = "'3.81"
$ws.Range("E37").Value = '  +2.91%  '
$ws.Range("E38").Value = '  -1.89%  '
$ws.Range("E39").Value = '  -2.26%  '
$ws.Range("E40").Value = '  -2.48%  '
$ws.Range("D41").Value = "'5.82"
$ws.Range("E41").Value = '  -2.52%  '
$ws.Range("E43").Value = '  -1.59%  '
$ws.Range("E44").Value = '  +0.00%  '
$ws.Range("D45").Value = "'8.69"
$ws.Range("E45").Value = '  -0.94%  '
$ws.Range("E46").Value = '  -2.42%  '
$ws.Range("D47").Value = "'45.88"
$ws.Range("E47").Value = '  -2.74%  '
$ws.Range("D48").Value = "'397.57"
$ws.Range("E48").Value = '  -5.04%  '
$ws.Range("E49").Value = '  -8.65%  '
$ws.Range("D50").Value = "'138.99"
$ws.Range("E50").Value = '  -1.93%  '
$ws.Range("E51").Value = '  -2.36%  '
